$wb = $excel.ActiveWorkbook

# Rename the "fielddomain" sheet to "field-domain"
$ws = $wb.Worksheets.Item("fielddomain")
$ws.Name = "field-domain"

# Select the renamed sheet and a cell on it, making it the active tab
$ws.Activate()
$ws.Range("F5").Select()
